$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the end time of the existing (split) session recorded on row 92 -
# it actually ended earlier than first logged.
$ws.Range("E92").Value = 0.4236111111111111

# Insert a new row to record the second work session of that same day
# (this pushes the blank separator row / summary rows down by one).
$ws.Rows.Item(93).Insert()

$ws.Range("A93").Value = 2014
$ws.Range("B93").Value = 3
$ws.Range("C93").Value = 25
$ws.Range("D93").Value = 0.47916666666666669
$ws.Range("E93").Value = 0.54166666666666663
$ws.Range("F93").Formula = "=(E93-D93)*24*60"
$ws.Range("G93").Formula = "=F93/60"

# Move the active selection to the (now empty) row below the new entry,
# matching where Excel would have left the cursor after the insert.
$ws.Range("A94").Select() | Out-Null
